# "correct inputs, add AR28"
# Personnel sheet: replace Kevin Cahill (technician) with Zoe Sandwith (creator),
# dropping the now-unused email address, and update the sheet's selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")
$ws.Activate()

$ws.Range("A6").Value = "Zoe"
$ws.Range("C6").Value = "Sandwith"
$ws.Range("E6").ClearContents()
$ws.Range("G6").Value = "creator"

$ws.Range("G7").Select()
